$wb = $excel.ActiveWorkbook

# --- Sheet: Energies ---
$ws1 = $wb.Worksheets.Item("Energies")
$ws1.Range("B2").Value = 13.693
$ws1.Range("C2").Value = 6.032
$ws1.Range("D2").Value = 46.524
$ws1.Range("D5").Value = 10.876
$ws1.Range("B6").Value = 11.915
$ws1.Range("C6").Value = 0.07000000000000001
$ws1.Range("D6").Value = 0.011

# --- Sheet: Partition Functions ---
$ws2 = $wb.Worksheets.Item("Partition Functions")
$ws2.Range("B4").Value = 53.1545
$ws2.Range("C4").Value = 1.72554
$ws2.Range("D4").Value = 3.973203
$ws2.Range("B5").Value = 0.497165
$ws2.Range("C5").Value = -0.303499
$ws2.Range("D5").Value = -0.698833
$ws2.Range("B6").Value = 268073000
$ws2.Range("C6").Value = 8.428253
$ws2.Range("D6").Value = 19.406769
$ws2.Range("B7").Value = [double]"1.85581e-09"
$ws2.Range("C7").Value = -8.731467
$ws2.Range("D7").Value = -20.104945
$ws2.Range("B8").Value = 1.00066
$ws2.Range("C8").Value = 0.000285
$ws2.Range("D8").Value = 0.000657

# --- Sheet: Other ---
$ws3 = $wb.Worksheets.Item("Other")
$ws3.Range("A2").Value = -55.8846012853
$ws3.Range("A3").Value = -55.8854240101
$ws3.Range("A4").Value = -55.8854468845
$ws3.Range("A5").Value = -55.8854468845

$ws3.Range("A14").Value = 2.56951
$ws3.Range("B14").Value = 4.71312
$ws3.Range("C14").Value = 7.28263

$ws3.Range("A23").Value = 33.70824
$ws3.Range("B23").Value = 18.37718
$ws3.Range("C23").Value = 11.8932

$ws3.Range("A26").Value = 702.36667
$ws3.Range("B26").Value = 382.91883
$ws3.Range("C26").Value = 247.81445

$ws3.Range("A29").Value = 49841

$ws3.Range("A32").Value = 11.91228

$ws3.Range("A35").Value = 2184.78
$ws3.Range("B35").Value = 4828.4
$ws3.Range("C35").Value = 4975.79

$ws3.Range("A38").Value = 0.018983

$ws3.Range("A41").Value = 0.021821

$ws3.Range("A44").Value = 0.022765

$ws3.Range("A47").Value = 0.00066

$ws3.Range("A50").Value = 1518.503
$ws3.Range("B50").Value = 3355.9093
$ws3.Range("C50").Value = 3458.3495

$ws3.Range("A53").Value = 1.0917
$ws3.Range("B53").Value = 1.0498
$ws3.Range("C53").Value = 1.0902

$ws3.Range("A56").Value = 1.4831
$ws3.Range("B56").Value = 6.9659
$ws3.Range("C56").Value = 7.6825

$ws3.Range("A59").Value = 23.3897
$ws3.Range("B59").Value = 16.2446
$ws3.Range("C59").Value = 2.4019
